$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
try {
  $s.ApplyThemeColorScheme("Office")
  Write-Host "ok"
} catch { Write-Host "err: $_" }
